$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.246358513832092
$ws.Range("B1").Value = 2.369637489318848
$ws.Range("C1").Value = 3.309111356735229
$ws.Range("D1").Value = 1.653137564659119
$ws.Range("E1").Value = 1.1922367811203
